$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2:B2) currently holds text "5000.0" values; convert them to real numbers (5000)
$ws.Range("A2").Value = 5000
$ws.Range("B2").Value = 5000

# Insert a new row of data (row 3) as requested by the commit "Inserindo dados na lista"
# A3 / B3 should hold the text value "4000.0" (kept as text, same style as other data cells)
$ws.Range("A3:B3").NumberFormat = "@"
$ws.Range("A3").Value = "4000.0"
$ws.Range("B3").Value = "4000.0"
$ws.Range("A3:B3").Style = "Normal"
